# Update crypto price/volume figures per the Oct 25 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.602.24"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "1.789.38"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Formula = "'224.85"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Formula = "'0.561"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Formula = "'32.55"
$ws.Range("E8").Value = "  +5.49%  "
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("D10").Formula = "'0.0669"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Formula = "'11.05"
$ws.Range("E13").Value = "  +10.09%  "
$ws.Range("D14").Value = "1.800.73"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Formula = "'0.633"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.623.79"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").Formula = "'68.85"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Formula = "'254.10"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").Value = "0.0₃0766"
$ws.Range("E20").Value = "  +3.19%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").Formula = "'159.51"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("E26").Value = "  -0.75%  "
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").Formula = "'0.0517"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("E34").Value = "  +3.04%  "
$ws.Range("D35").Value = "1.442.62"
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Formula = "'0.0189"
$ws.Range("E37").Value = "  +2.01%  "
$ws.Range("D38").Formula = "'0.628"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("D39").Formula = "'83.14"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("E40").Value = "  +4.22%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Formula = "'0.0504"
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").Value = "1.945.86"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").Formula = "'11.99"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").Formula = "'103.06"
$ws.Range("E50").Value = "  +5.21%  "
$ws.Range("E51").Value = "  +6.05%  "
